# Add download suppliers orders:
#  - remove the "Location" column (old column B)
#  - remove the old row 3 (Тунец / Холодильник 2 / кг / 3)
#  - replace the remaining data row with the new supplier order entry
#    (Курица / л / 12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Location" column entirely; remaining columns shift left so the
# old "Measurement"/"Quantity" columns become B/C.
$ws.Columns.Item(2).Delete()

# Drop the old last row (was row 3: Тунец / Холодильник 2 / кг / 3).
$ws.Rows.Item(3).Delete()

# Overwrite the remaining data row with the new order.
$ws.Range("A2").Value = "Курица"
$ws.Range("B2").Value = "л"
$ws.Range("C2").Value = 12
